$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.119.38'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '2.925.62'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.39'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.31'
$ws.Range("E6").Value = '  -1.29%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  +0.95%  '
$ws.Range("D9").Value = '2.926.47'
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.83'
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.144'
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.443'
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("E13").Value = '  +1.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.75'
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '3.415.23'
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '61.081.34'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.73'
$ws.Range("E18").Value = '  -1.93%  '
$ws.Range("D19").Value = '2.925.15'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '431.39'
$ws.Range("E20").Value = '  +0.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.50'
$ws.Range("E21").Value = '  -2.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.683'
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.10'
$ws.Range("E23").Value = '  -1.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.92'
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.83'
$ws.Range("E25").Value = '  -0.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("E26").Value = '  +1.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.17'
$ws.Range("E27").Value = '  +2.00%  '
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.34'
$ws.Range("E29").Value = '  +6.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.14'
$ws.Range("E32").Value = '  -2.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.59'
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("E34").Value = '  +0.79%  '
$ws.Range("D35").Value = '0.0₃0866'
$ws.Range("E35").Value = '  +3.97%  '
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.11'
$ws.Range("E37").Value = '  +3.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.63'
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.99'
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.01'
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("E41").Value = '  -1.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.63'
$ws.Range("E42").Value = '  -1.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.289'
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.89'
$ws.Range("E44").Value = '  -4.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '383.09'
$ws.Range("E45").Value = '  +2.83%  '
$ws.Range("E46").Value = '  +0.53%  '
$ws.Range("D47").Value = '2.708.16'
$ws.Range("E47").Value = '  +1.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.88'
$ws.Range("E48").Value = '  -2.55%  '
$ws.Range("E50").Value = '  -4.07%  '
$ws.Range("E51").Value = '  +0.26%  '

# Restore default style on cells where we forced text NumberFormat
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Style = "Normal"
